# Update the "Who?" (assignee) column on the Planning sheet:
#  - Row 13 "Determine keypoints in images (and save to db)": add Gillis
#  - Row 14 "Determine feature vectors in images (and save to db)": fix name typo,
#    now matches row 13's assignees
#  - Row 15 "Predict room with test images": add Thomas (was unassigned)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")

$ws.Range("G13").Value = "Gillis, Pieter-Jan, Thomas"
$ws.Range("G14").Value = "Gillis, Pieter-Jan, Thomas"
$ws.Range("G15").Value = "Thomas"
